$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").Value = 9.5695481095730948
$ws.Range("C2").Value = -0.89239373006500067
$ws.Range("D2").Value = 0.19794520859791626
$ws.Range("E2").Value = -0.33222657568753378

$ws.Range("B3").Value = 2.7472038679327975
$ws.Range("C3").Value = 4.0238409481603599
$ws.Range("D3").Value = 2.2223642054436512
$ws.Range("E3").Value = -1.6192884509807115

$ws.Range("B1:E3").Select()
